$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (prices, volumes, and re-ranked rows 9-27)

$ws.Range("D2").Value = "'256.72"
$ws.Range("E2").Value = "'-0.71%"

$ws.Range("D3").Value = "'26.98"
$ws.Range("E3").Value = "'0.17%"

$ws.Range("D4").Value = "'4.389"
$ws.Range("E4").Value = "'-5.90%"

$ws.Range("D5").Value = "'0.05892"

$ws.Range("D6").Value = "'6.637"
$ws.Range("E6").Value = "'-0.72%"

$ws.Range("E7").Value = "'-2.86%"

$ws.Range("D8").Value = "'0.9386"
$ws.Range("E8").Value = "'-1.74%"

$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01035"
$ws.Range("E9").Value = "'1,600.30%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1383"
$ws.Range("E10").Value = "'-2.33%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.04814"
$ws.Range("E11").Value = "'33.43%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07086"
$ws.Range("E12").Value = "'-1.41%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03076"
$ws.Range("E13").Value = "'-1.95%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09106"
$ws.Range("E14").Value = "'-1.38%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001532"
$ws.Range("E15").Value = "'-0.72%"

$ws.Range("D16").Value = "'0.006173"
$ws.Range("E16").Value = "'3.70%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.492"
$ws.Range("E17").Value = "'0.14%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.168"
$ws.Range("E18").Value = "'-1.72%"

$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.226"
$ws.Range("E19").Value = "'0.30%"

$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3049"
$ws.Range("E20").Value = "'-2.78%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1270"
$ws.Range("E21").Value = "'-1.59%"

$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'3.921"
$ws.Range("E22").Value = "'10.99%"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04279"
$ws.Range("E23").Value = "'1.36%"

$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'-0.09%"

$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004282"
$ws.Range("E25").Value = "'-5.16%"

$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.05%"

$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").Value = "'0.0001524"
$ws.Range("E27").Value = "'2.08%"

$ws.Range("D40").Value = "'0.03815"
$ws.Range("E40").Value = "'-0.89%"

$ws.Range("D41").Value = "'0.006201"
$ws.Range("E41").Value = "'2.87%"

$ws.Range("D42").Value = "'0.1100"
$ws.Range("E42").Value = "'-0.38%"

$ws.Range("D43").Value = "'0.002200"
$ws.Range("E43").Value = "'0.04%"

$ws.Range("D44").Value = "'0.01398"
$ws.Range("E44").Value = "'26.16%"

$ws.Range("D45").Value = "'0.00005373"
$ws.Range("E45").Value = "'-2.21%"

$ws.Range("E46").Value = "'0.04%"

$ws.Range("D47").Value = "'0.06589"
$ws.Range("E47").Value = "'-22.92%"

$ws.Range("D48").Value = "'0.2517"
$ws.Range("E48").Value = "'11,738.47%"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.04%"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.04%"
